# Apply attendance-marking updates to Sheet1.
# For each listed cell, set its value from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @(
    "G3", "H3",
    "D4", "E4",
    "D5", "E5",
    "D6", "E6",
    "H7",
    "H8",
    "H9",
    "H10",
    "H11",
    "D12", "E12",
    "H13",
    "D14", "E14",
    "D15", "E15",
    "H16",
    "D17", "E17",
    "H18"
)

foreach ($cellAddr in $cells) {
    $ws.Range($cellAddr).Value = 1
}
